# Applies the "Fixed bugs and added corner cases" edit to test.xlsx
# Sheet1: fill in new corner-case rows for the Analog/Digital/Measurements
#         sections and add a brand-new "{Utility ... }" section.
# Sheet2: populate what was an empty scratch sheet with a second worked
#         example (Configure / configure corner cases).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------
# Sheet1 edits
# ---------------------------------------------------------------

# New comment line under the header block
$ws1.Range("A6").Value = "This is a comment"

# New "analog" comment above the Configure Frequency/Amplitude rows
$ws1.Range("C10").Value = "#This is analog"

# Row 11 (Configure Frequency): tweak the input spec, replace the SCPI
# command with a comment, and drop the trailing comment cell entirely
$ws1.Range("D11").Value = "in: Frequency:double"
$ws1.Range("E11").Value = "#This is a comment"
$ws1.Range("F11").ClearContents()

# Row 12 (Configure Amplitude): drop the leading "# " from the comment
$ws1.Range("F12").Value = "Sets the amplitude of the analog measurement."

# Row 16 (Configure dBM): add input spec + SCPI command
$ws1.Range("D16").Value = "in:dBM:integer"
$ws1.Range("E16").Value = "MIT:TEST:DBM <value>"

# Row 17 (Configure Resolution): add input spec + SCPI command
$ws1.Range("D17").Value = "in:Resoluttion:Double"
$ws1.Range("E17").Value = "MIT:TEST:RESO <value>"

# Measurement rows 25-28: add output specs
$ws1.Range("D25").Value = "out:Measurement:Ring:val1,val2,"
$ws1.Range("D26").Value = "out:Frequency:int"
$ws1.Range("D27").Value = "out:Last Data:string"
$ws1.Range("D28").Value = "none"

# New "{Utility ... }" section (rows 33-35), mirroring the existing
# "{Configure" / "{Action-Status}" / "{Data" section headers: 14pt font,
# with only the word after "{" in bold.
$utilHeader = $ws1.Range("A33")
$utilHeader.Value = "{Utility"
$utilHeader.Font.Size = 14
$utilHeader.Characters(1, 1).Font.Size = 14
$utilHeader.Characters(2, 7).Font.Size = 14
$utilHeader.Characters(2, 7).Font.Bold = $true
$ws1.Rows.Item(33).RowHeight = 18.75

$ws1.Range("B34").Value = "Reset"
$ws1.Range("C34").Value = "none"
$ws1.Range("D34").Value = "*RST"

$ws1.Range("A35").Value = "}"

# ---------------------------------------------------------------
# Sheet2 edits - second worked example, previously blank
# ---------------------------------------------------------------

$ws1.Range("A1").Copy() | Out-Null
$ws2.Range("A2").PasteSpecial(-4104) | Out-Null
$ws1.Range("A2").Copy() | Out-Null
$ws2.Range("A3").PasteSpecial(-4104) | Out-Null
$ws1.Range("A3").Copy() | Out-Null
$ws2.Range("A4").PasteSpecial(-4104) | Out-Null
$ws1.Range("A4").Copy() | Out-Null
$ws2.Range("A5").PasteSpecial(-4104) | Out-Null

$ws2.Range("A7").Value = "{Configure"
$ws2.Range("B8").Value = "Im a Vi"
$ws2.Range("A9").Value = "}"

$ws2.Range("A11").Value = "{configure"

$ws2.Columns.Item(1).ColumnWidth = 27.92
$ws2.Columns.Item(2).ColumnWidth = 24.42
$ws2.Columns.Item(3).ColumnWidth = 20.92
$ws2.Columns.Item(4).ColumnWidth = 23.09
$ws2.Columns.Item(5).ColumnWidth = 44.09
$ws2.Columns.Item(6).ColumnWidth = 9.09
$ws2.Columns.Item(7).ColumnWidth = 35.09

# ---------------------------------------------------------------
# Selection / view state
# ---------------------------------------------------------------

$ws2.Range("A13").Select()

$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 13
$ws1.Range("E28").Select()
